$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functional Testing")

# Fix the typo "ransom.html" -> "random.html" in the two cells that reference this text
$oldText = "When `"We've got to the end of that little selection`" is displayed the `"Show me a random selection`" button links to ransom.html"
$newText = "When `"We've got to the end of that little selection`" is displayed the `"Show me a random selection`" button links to random.html"

$ws.Range("A40").Value = $newText
$ws.Range("A63").Value = $newText

# Update the view: scroll so row 52 is at top-left, and select B71
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B71").Select()
